$p = $ppt.ActivePresentation

# --- Update the cached "last updated" date placeholder text on the slide
#     master and on every slide layout (6/16/2025 -> 6/18/2025). ---
function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = "6/18/2025"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li)
}

# --- Fix the "res_types" label (typo) on slide 1 and shrink its textbox to
#     match the now-shorter auto-fit text width. ---
$slide = $p.Slides.Item(1)
$label = $slide.Shapes.Item("TextBox 51")
$label.TextFrame.TextRange.Text = "res_type"
$label.Width = 67.048465
